$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.061.52"

$ws.Range("D3").Value = "1.650.84"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5220"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06276"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("E10").Value = "  -3.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07735"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.643.77"
$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.465"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5440"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").Value = "0.0₅8098"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "26.076.70"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.575"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.992"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "138.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.261"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.399"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05939"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.03%  "

$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.498"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.228"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.537"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9473"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.72%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.753"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5729"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.64%  "

$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.863"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8464"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.66%  "

$ws.Range("D43").Value = "1.003.08"
$ws.Range("E43").Value = "  -4.49%  "

$ws.Range("D44").Value = "1.794.55"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("E46").Value = "  -2.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4300"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.903"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("E51").Value = "  -0.71%  "
